$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$entries = @(
    @{ Cell = 'D2'; Value = '95.594.21' },
    @{ Cell = 'E2'; Value = '  -1.80%  ' },
    @{ Cell = 'D3'; Value = '3.611.19' },
    @{ Cell = 'E3'; Value = '  -2.45%  ' },
    @{ Cell = 'E4'; Value = '  +22.73%  ' },
    @{ Cell = 'E5'; Value = '  +0.01%  ' },
    @{ Cell = 'D6'; Value = '222.67' },
    @{ Cell = 'E6'; Value = '  -5.93%  ' },
    @{ Cell = 'D7'; Value = '640.05' },
    @{ Cell = 'E7'; Value = '  -2.45%  ' },
    @{ Cell = 'D8'; Value = '0.421' },
    @{ Cell = 'E8'; Value = '  -4.54%  ' },
    @{ Cell = 'E9'; Value = '  +2.09%  ' },
    @{ Cell = 'D10'; Value = '0.999' },
    @{ Cell = 'E10'; Value = '  +0.01%  ' },
    @{ Cell = 'D11'; Value = '3.605.45' },
    @{ Cell = 'E11'; Value = '  -2.62%  ' },
    @{ Cell = 'D12'; Value = '49.66' },
    @{ Cell = 'E12'; Value = '  +10.08%  ' },
    @{ Cell = 'D13'; Value = '0.216' },
    @{ Cell = 'E13'; Value = '  +4.05%  ' },
    @{ Cell = 'D14'; Value = '0.0000290' },
    @{ Cell = 'E14'; Value = '  -6.74%  ' },
    @{ Cell = 'D15'; Value = '6.49' },
    @{ Cell = 'E15'; Value = '  -5.02%  ' },
    @{ Cell = 'D16'; Value = '25.38' },
    @{ Cell = 'E16'; Value = '  +33.99%  ' },
    @{ Cell = 'D17'; Value = '4.281.54' },
    @{ Cell = 'D18'; Value = '95.329.26' },
    @{ Cell = 'E18'; Value = '  -1.67%  ' },
    @{ Cell = 'D19'; Value = '9.13' },
    @{ Cell = 'E19'; Value = '  +2.61%  ' },
    @{ Cell = 'E20'; Value = '  +4.99%  ' },
    @{ Cell = 'D21'; Value = '3.595.17' },
    @{ Cell = 'E21'; Value = '  -3.72%  ' },
    @{ Cell = 'D22'; Value = '0.280' },
    @{ Cell = 'E22'; Value = '  +36.98%  ' },
    @{ Cell = 'D23'; Value = '0.530' },
    @{ Cell = 'E23'; Value = '  -2.20%  ' },
    @{ Cell = 'D24'; Value = '136.20' },
    @{ Cell = 'E24'; Value = '  +15.60%  ' },
    @{ Cell = 'D25'; Value = '529.52' },
    @{ Cell = 'E26'; Value = '  -5.25%  ' },
    @{ Cell = 'E27'; Value = '  -9.12%  ' },
    @{ Cell = 'D28'; Value = '6.89' },
    @{ Cell = 'E28'; Value = '  -0.36%  ' },
    @{ Cell = 'D29'; Value = '3.773.96' },
    @{ Cell = 'E29'; Value = '  -3.23%  ' },
    @{ Cell = 'D30'; Value = '13.00' },
    @{ Cell = 'E30'; Value = '  -3.20%  ' },
    @{ Cell = 'D31'; Value = '13.22' },
    @{ Cell = 'E31'; Value = '  +4.05%  ' },
    @{ Cell = 'D32'; Value = '3.15' },
    @{ Cell = 'E32'; Value = '  +3.94%  ' },
    @{ Cell = 'E33'; Value = '  +0.08%  ' },
    @{ Cell = 'E34'; Value = '  +7.08%  ' },
    @{ Cell = 'E35'; Value = '  +2.61%  ' },
    @{ Cell = 'D36'; Value = '33.69' },
    @{ Cell = 'E36'; Value = '  +1.88%  ' },
    @{ Cell = 'D37'; Value = '0.184' },
    @{ Cell = 'E37'; Value = '  -2.77%  ' },
    @{ Cell = 'E38'; Value = '  +0.19%  ' },
    @{ Cell = 'E39'; Value = '  -0.01%  ' },
    @{ Cell = 'B40'; Value = 'RenderToken' },
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render' },
    @{ Cell = 'D40'; Value = '8.45' },
    @{ Cell = 'E40'; Value = '  -3.25%  ' },
    @{ Cell = 'B41'; Value = 'Filecoin' },
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil' },
    @{ Cell = 'D41'; Value = '7.21' },
    @{ Cell = 'E41'; Value = '  +5.16%  ' },
    @{ Cell = 'B42'; Value = 'Bittensor' },
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao' },
    @{ Cell = 'D42'; Value = '591.07' },
    @{ Cell = 'E42'; Value = '  -6.66%  ' },
    @{ Cell = 'D43'; Value = '0.0535' },
    @{ Cell = 'E43'; Value = '  +17.36%  ' },
    @{ Cell = 'B44'; Value = 'ARBITRUM' },
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb' },
    @{ Cell = 'D44'; Value = '1.04' },
    @{ Cell = 'E44'; Value = '  +7.88%  ' },
    @{ Cell = 'B45'; Value = 'Algorand' },
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo' },
    @{ Cell = 'D45'; Value = '0.502' },
    @{ Cell = 'E45'; Value = '  -1.08%  ' },
    @{ Cell = 'D46'; Value = '41.16' },
    @{ Cell = 'D47'; Value = '0.160' },
    @{ Cell = 'E47'; Value = '  -3.57%  ' },
    @{ Cell = 'D48'; Value = '1.98' },
    @{ Cell = 'E48'; Value = '  -1.47%  ' },
    @{ Cell = 'D49'; Value = '9.25' },
    @{ Cell = 'E49'; Value = '  +4.89%  ' },
    @{ Cell = 'D50'; Value = '232.71' },
    @{ Cell = 'E50'; Value = '  +13.03%  ' },
    @{ Cell = 'D51'; Value = '2.33' },
    @{ Cell = 'E51'; Value = '  -2.43%  ' }
)

foreach ($entry in $entries) {
    $cell = $ws.Range($entry.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $entry.Value
}
